$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timestamp in A1
$ws.Range("A1").Value = 'Datos actualizados a 5 de Septiembre de 2020 a las 18:04'

# Row updates: country name (col A) plus stats columns B-H
$rows = @(
  @{Row=4; Vals=@('Estados Unidos', 6400414, 11357, 3637002, 2571104, 0, 197, 192308)},
  @{Row=5; Vals=@('India', 4092550, 72311, 3162305, 859726, 0, 884, 70519)},
  @{Row=6; Vals=@('Brasil', 4091801, 0, 3278243, 687974, 0, 0, 125584)},
  @{Row=14; Vals=@('Chile', 420434, 1969, 392967, 15916, 0, 57, 11551)},
  @{Row=16; Vals=@('Reino Unido', 344164, 1813, 0, 0, 0, 12, 41549)},
  @{Row=22; Vals=@('Italia', 276338, 1695, 209610, 31194, 0, 16, 35534)},
  @{Row=28; Vals=@('Canada', 131467, 343, 116118, 6206, 0, 2, 9143)},
  @{Row=35; Vals=@('Republica Dominicana', 98776, 864, 71901, 25035, 0, 19, 1840)},
  @{Row=52; Vals=@('Singapur', 56982, 34, 56267, 688, 0, 0, 27)},
  @{Row=63; Vals=@('Suiza', 43957, 425, 37100, 4844, 0, 0, 2013)},
  @{Row=64; Vals=@('Uzbekistan', 43293, 295, 40774, 2176, 0, 5, 343)},
  @{Row=65; Vals=@('Moldavia', 39473, 567, 27799, 10611, 0, 16, 1063)},
  @{Row=94; Vals=@('Albania', 10102, 135, 5976, 3814, 0, 6, 312)},
  @{Row=106; Vals=@('Luxemburgo', 6854, 0, 6126, 604, 0, 0, 124)},
  @{Row=119; Vals=@('Mozambique', 4341, 76, 2579, 1736, 0, 0, 26)},
  @{Row=120; Vals=@('Ruanda', 4304, 0, 2191, 2095, 0, 0, 18)},
  @{Row=121; Vals=@('Cuba', 4266, 0, 3487, 679, 0, 0, 100)},
  @{Row=130; Vals=@('Sri Lanka', 3118, 3, 2918, 188, 0, 0, 12)},
  @{Row=140; Vals=@('Jordania', 2353, 52, 1700, 637, 0, 0, 16)},
  @{Row=143; Vals=@('Trinidad yTobago', 2142, 102, 717, 1393, 0, 1, 32)},
  @{Row=144; Vals=@('Islandia', 2136, 1, 2038, 88, 0, 0, 10)},
  @{Row=145; Vals=@('Reunion', 2115, 113, 1313, 791, 0, 1, 11)},
  @{Row=146; Vals=@('Sierra Leona', 2041, 0, 1602, 368, 0, 0, 71)},
  @{Row=202; Vals=@('Fiyi', 31, 2, 24, 5, 0, 0, 2)},
  @{Row=207; Vals=@('Dominica', 22, 2, 18, 4, 0, 0, 0)},
  @{Row=208; Vals=@('Laos', 22, 0, 21, 1, 0, 0, 0)}
)

foreach ($r in $rows) {
  $rowNum = $r.Row
  $vals = $r.Vals
  for ($c = 0; $c -lt $vals.Length; $c++) {
    $ws.Cells.Item($rowNum, $c + 1).Value = $vals[$c]
  }
}

Write-Host "Applied country/provincias update"
